$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.447.30"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.426.91"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.25"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.35"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.78"
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.148"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("E13").Value = "  +4.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.859.45"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.364.99"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.436.68"
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.42"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.27"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.36"
$ws.Range("E23").Value = "  +4.40%  "
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.62"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.36"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0788"
$ws.Range("E28").Value = "  +6.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.80"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").Value = "  +2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.42"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.78"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.32"
$ws.Range("E35").Value = "  +6.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.24"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.92"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.419"
$ws.Range("E40").Value = "  +10.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "319.38"
$ws.Range("E41").Value = "  +9.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.73"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "142.65"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0525"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0959"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.77"
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.414"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.00"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("E51").Value = "  -0.21%  "
